$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Players")

# Column G (status) width: 17 -> 18
$ws.Columns.Item(7).ColumnWidth = 17.15

# Row 2
$ws.Cells.Item(2, 7).Value = "16:31 - 2nd Half"
$ws.Cells.Item(2, 8).Value = 6
$ws.Cells.Item(2, 15).Value = 3
$ws.Cells.Item(2, 16).Value = 20
$ws.Cells.Item(2, 18).Value = 8
$ws.Cells.Item(2, 20).Value = 4

# Row 3
$ws.Cells.Item(3, 7).Value = "16:31 - 2nd Half"
$ws.Cells.Item(3, 16).Value = 22

# Row 6
$ws.Cells.Item(6, 7).Value = "Final"

# Row 7
$ws.Cells.Item(7, 7).Value = "16:31 - 2nd Half"
$ws.Cells.Item(7, 8).Value = 8
$ws.Cells.Item(7, 12).Value = 3
$ws.Cells.Item(7, 16).Value = 20
$ws.Cells.Item(7, 18).Value = 10
$ws.Cells.Item(7, 20).Value = 4

# Row 9
$ws.Cells.Item(9, 7).Value = "16:31 - 2nd Half"
$ws.Cells.Item(9, 15).Value = 3

# Row 11
$ws.Cells.Item(11, 7).Value = "Final"

# Row 12
$ws.Cells.Item(12, 7).Value = "16:31 - 2nd Half"
$ws.Cells.Item(12, 10).Value = 1
$ws.Cells.Item(12, 15).Value = 1
$ws.Cells.Item(12, 16).Value = 19
$ws.Cells.Item(12, 18).Value = 4

# Row 13
$ws.Cells.Item(13, 7).Value = "Final"

# Row 18
$ws.Cells.Item(18, 7).Value = "16:31 - 2nd Half"
$ws.Cells.Item(18, 8).Value = 7
$ws.Cells.Item(18, 9).Value = 7
$ws.Cells.Item(18, 10).Value = 6
$ws.Cells.Item(18, 15).Value = 2
$ws.Cells.Item(18, 16).Value = 20
$ws.Cells.Item(18, 21).Value = 3
$ws.Cells.Item(18, 22).Value = 4

# Row 19
$ws.Cells.Item(19, 7).Value = "Final"

# Row 21
$ws.Cells.Item(21, 7).Value = "16:31 - 2nd Half"
$ws.Cells.Item(21, 16).Value = 14

# Row 23
$ws.Cells.Item(23, 7).Value = "Final"
$ws.Cells.Item(23, 8).Value = -1
$ws.Cells.Item(23, 9).Value = 13
$ws.Cells.Item(23, 16).Value = 37
$ws.Cells.Item(23, 18).Value = 20
$ws.Cells.Item(23, 20).Value = 12
$ws.Cells.Item(23, 21).Value = 3
$ws.Cells.Item(23, 22).Value = 5

# Row 24
$ws.Cells.Item(24, 7).Value = "Final"

# Row 25
$ws.Cells.Item(25, 7).Value = "16:31 - 2nd Half"
$ws.Cells.Item(25, 8).Value = 15
$ws.Cells.Item(25, 9).Value = 14
$ws.Cells.Item(25, 14).Value = 3
$ws.Cells.Item(25, 16).Value = 21
$ws.Cells.Item(25, 17).Value = 5
$ws.Cells.Item(25, 18).Value = 5
$ws.Cells.Item(25, 19).Value = 2
$ws.Cells.Item(25, 20).Value = 2

# Row 26
$ws.Cells.Item(26, 7).Value = "Final"
$ws.Cells.Item(26, 8).Value = 26
$ws.Cells.Item(26, 9).Value = 21
$ws.Cells.Item(26, 10).Value = 4
$ws.Cells.Item(26, 21).Value = 8
$ws.Cells.Item(26, 22).Value = 10

# Row 27
$ws.Cells.Item(27, 7).Value = "Final"

# Row 28
$ws.Cells.Item(28, 7).Value = "16:31 - 2nd Half"
$ws.Cells.Item(28, 8).Value = 8
$ws.Cells.Item(28, 10).Value = 6
$ws.Cells.Item(28, 13).Value = 4
$ws.Cells.Item(28, 14).Value = 2
$ws.Cells.Item(28, 16).Value = 16
$ws.Cells.Item(28, 22).Value = 2

# Row 29
$ws.Cells.Item(29, 7).Value = "Final"
$ws.Cells.Item(29, 10).Value = 3
$ws.Cells.Item(29, 14).Value = 1
$ws.Cells.Item(29, 16).Value = 23

# Row 30
$ws.Cells.Item(30, 7).Value = "Final"
$ws.Cells.Item(30, 8).Value = 16
$ws.Cells.Item(30, 12).Value = 2
$ws.Cells.Item(30, 16).Value = 28

# Row 32
$ws.Cells.Item(32, 7).Value = "Final"

# Row 34
$ws.Cells.Item(34, 7).Value = "Final"

# Row 35
$ws.Cells.Item(35, 7).Value = "Final"

# Row 36
$ws.Cells.Item(36, 4).Value = "Bishop Boswell"
$ws.Cells.Item(36, 5).Value = "TENN"
$ws.Cells.Item(36, 6).Value = "TENN@VAN"
$ws.Cells.Item(36, 7).Value = "16:31 - 2nd Half"
$ws.Cells.Item(36, 9).Value = 5
$ws.Cells.Item(36, 10).Value = 2
$ws.Cells.Item(36, 11).Value = 3
$ws.Cells.Item(36, 12).Value = 2
$ws.Cells.Item(36, 13).Value = 0
$ws.Cells.Item(36, 14).Value = 1
$ws.Cells.Item(36, 15).Value = 2
$ws.Cells.Item(36, 16).Value = 16
$ws.Cells.Item(36, 18).Value = 4
$ws.Cells.Item(36, 19).Value = 0
$ws.Cells.Item(36, 20).Value = 0

# Row 37
$ws.Cells.Item(37, 4).Value = "Jordan Butler"
$ws.Cells.Item(37, 5).Value = "SC"
$ws.Cells.Item(37, 6).Value = "MSST@SC"
$ws.Cells.Item(37, 7).Value = "Final"
$ws.Cells.Item(37, 8).Value = 8
$ws.Cells.Item(37, 10).Value = 3
$ws.Cells.Item(37, 13).Value = 1
$ws.Cells.Item(37, 15).Value = 3
$ws.Cells.Item(37, 16).Value = 10
$ws.Cells.Item(37, 18).Value = 3
$ws.Cells.Item(37, 19).Value = 1
$ws.Cells.Item(37, 21).Value = 1
$ws.Cells.Item(37, 22).Value = 2

# Row 38
$ws.Cells.Item(38, 4).Value = "DeWayne Brown II"
$ws.Cells.Item(38, 5).Value = "TENN"
$ws.Cells.Item(38, 6).Value = "TENN@VAN"
$ws.Cells.Item(38, 7).Value = "16:31 - 2nd Half"
$ws.Cells.Item(38, 8).Value = 7
$ws.Cells.Item(38, 9).Value = 4
$ws.Cells.Item(38, 10).Value = 3
$ws.Cells.Item(38, 11).Value = 0
$ws.Cells.Item(38, 12).Value = 1
$ws.Cells.Item(38, 13).Value = 1
$ws.Cells.Item(38, 14).Value = 0
$ws.Cells.Item(38, 16).Value = 17
$ws.Cells.Item(38, 17).Value = 2
$ws.Cells.Item(38, 18).Value = 4
$ws.Cells.Item(38, 20).Value = 0
$ws.Cells.Item(38, 21).Value = 0
$ws.Cells.Item(38, 22).Value = 0

# Row 39
$ws.Cells.Item(39, 4).Value = "King Grace"
$ws.Cells.Item(39, 5).Value = "MSST"
$ws.Cells.Item(39, 6).Value = "MSST@SC"
$ws.Cells.Item(39, 7).Value = "Final"
$ws.Cells.Item(39, 8).Value = 7
$ws.Cells.Item(39, 9).Value = 6
$ws.Cells.Item(39, 10).Value = 2
$ws.Cells.Item(39, 11).Value = 0
$ws.Cells.Item(39, 12).Value = 0
$ws.Cells.Item(39, 14).Value = 0
$ws.Cells.Item(39, 16).Value = 15
$ws.Cells.Item(39, 18).Value = 3
$ws.Cells.Item(39, 20).Value = 1
$ws.Cells.Item(39, 21).Value = 2
$ws.Cells.Item(39, 22).Value = 2

# Row 40
$ws.Cells.Item(40, 4).Value = "Chandler Bing"
$ws.Cells.Item(40, 5).Value = "VAN"
$ws.Cells.Item(40, 6).Value = "TENN@VAN"
$ws.Cells.Item(40, 7).Value = "16:31 - 2nd Half"
$ws.Cells.Item(40, 8).Value = 6
$ws.Cells.Item(40, 9).Value = 6
$ws.Cells.Item(40, 11).Value = 0
$ws.Cells.Item(40, 12).Value = 0
$ws.Cells.Item(40, 14).Value = 1
$ws.Cells.Item(40, 15).Value = 1
$ws.Cells.Item(40, 16).Value = 12
$ws.Cells.Item(40, 17).Value = 2
$ws.Cells.Item(40, 19).Value = 2
$ws.Cells.Item(40, 20).Value = 2

# Row 41
$ws.Cells.Item(41, 4).Value = "Patton Pinkins"
$ws.Cells.Item(41, 5).Value = "MISS"
$ws.Cells.Item(41, 8).Value = 6
$ws.Cells.Item(41, 9).Value = 10
$ws.Cells.Item(41, 10).Value = 2
$ws.Cells.Item(41, 11).Value = 1
$ws.Cells.Item(41, 15).Value = 2
$ws.Cells.Item(41, 16).Value = 36
$ws.Cells.Item(41, 17).Value = 3
$ws.Cells.Item(41, 18).Value = 8
$ws.Cells.Item(41, 20).Value = 3
$ws.Cells.Item(41, 21).Value = 4
$ws.Cells.Item(41, 22).Value = 4

# Row 42
$ws.Cells.Item(42, 4).Value = "EJ Walker"
$ws.Cells.Item(42, 5).Value = "SC"
$ws.Cells.Item(42, 7).Value = "Final"
$ws.Cells.Item(42, 8).Value = 5
$ws.Cells.Item(42, 10).Value = 1
$ws.Cells.Item(42, 11).Value = 1
$ws.Cells.Item(42, 12).Value = 1
$ws.Cells.Item(42, 15).Value = 5
$ws.Cells.Item(42, 16).Value = 16
$ws.Cells.Item(42, 19).Value = 1
$ws.Cells.Item(42, 21).Value = 0
$ws.Cells.Item(42, 22).Value = 0

# Row 43
$ws.Cells.Item(43, 4).Value = "Isaiah Brown"
$ws.Cells.Item(43, 5).Value = "FLA"
$ws.Cells.Item(43, 8).Value = 4
$ws.Cells.Item(43, 9).Value = 5
$ws.Cells.Item(43, 10).Value = 3
$ws.Cells.Item(43, 11).Value = 0
$ws.Cells.Item(43, 12).Value = 0
$ws.Cells.Item(43, 13).Value = 0
$ws.Cells.Item(43, 14).Value = 2
$ws.Cells.Item(43, 15).Value = 3
$ws.Cells.Item(43, 16).Value = 13
$ws.Cells.Item(43, 17).Value = 2
$ws.Cells.Item(43, 18).Value = 3
$ws.Cells.Item(43, 21).Value = 1
$ws.Cells.Item(43, 22).Value = 2

# Row 44
$ws.Cells.Item(44, 4).Value = "Jamarion Davis-Fleming"
$ws.Cells.Item(44, 5).Value = "MSST"
$ws.Cells.Item(44, 6).Value = "MSST@SC"
$ws.Cells.Item(44, 7).Value = "Final"
$ws.Cells.Item(44, 8).Value = 4
$ws.Cells.Item(44, 9).Value = 5
$ws.Cells.Item(44, 10).Value = 1
$ws.Cells.Item(44, 11).Value = 2
$ws.Cells.Item(44, 12).Value = 0
$ws.Cells.Item(44, 13).Value = 0
$ws.Cells.Item(44, 14).Value = 1
$ws.Cells.Item(44, 15).Value = 3
$ws.Cells.Item(44, 16).Value = 17
$ws.Cells.Item(44, 17).Value = 1
$ws.Cells.Item(44, 21).Value = 3
$ws.Cells.Item(44, 22).Value = 6

# Row 45
$ws.Cells.Item(45, 4).Value = "Corey Chest"
$ws.Cells.Item(45, 5).Value = "MISS"
$ws.Cells.Item(45, 6).Value = "FLA@MISS"
$ws.Cells.Item(45, 7).Value = "Final"
$ws.Cells.Item(45, 9).Value = 0
$ws.Cells.Item(45, 10).Value = 1
$ws.Cells.Item(45, 11).Value = 1
$ws.Cells.Item(45, 12).Value = 1
$ws.Cells.Item(45, 13).Value = 2
$ws.Cells.Item(45, 14).Value = 0
$ws.Cells.Item(45, 15).Value = 1
$ws.Cells.Item(45, 16).Value = 12
$ws.Cells.Item(45, 17).Value = 0
$ws.Cells.Item(45, 18).Value = 2
$ws.Cells.Item(45, 21).Value = 0
$ws.Cells.Item(45, 22).Value = 0

# Row 46
$ws.Cells.Item(46, 7).Value = "16:31 - 2nd Half"

# Row 50
$ws.Cells.Item(50, 7).Value = "16:31 - 2nd Half"
$ws.Cells.Item(50, 15).Value = 3
$ws.Cells.Item(50, 16).Value = 8

# Row 52
$ws.Cells.Item(52, 7).Value = "Final"
$ws.Cells.Item(52, 16).Value = 17

# Row 54
$ws.Cells.Item(54, 7).Value = "16:31 - 2nd Half"

# Row 55
$ws.Cells.Item(55, 7).Value = "16:31 - 2nd Half"

# Row 58
$ws.Cells.Item(58, 7).Value = "Final"

# Row 60
$ws.Cells.Item(60, 7).Value = "16:31 - 2nd Half"

# Row 63
$ws.Cells.Item(63, 7).Value = "Final"

# OwnerTotals sheet updates
$ws2 = $wb.Worksheets.Item("OwnerTotals")
$ws2.Cells.Item(2, 2).Value = 61
$ws2.Cells.Item(4, 2).Value = 15
$ws2.Cells.Item(5, 2).Value = 8
$ws2.Cells.Item(6, 2).Value = 6
